$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a TEXT cell even when $value looks numeric
    # (e.g. "1.00", "0.999"), matching the source inlineStr cells, then
    # restore the default "Normal" style so no stray number-format is left
    # attached to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '70.302.60'
$ws.Range("E2").Value = '  -0.52%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.587.12'
$ws.Range("E3").Value = '  -1.29%  '

# Row 4
Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
Set-TextValue $ws.Range("D5") '580.35'
$ws.Range("E5").Value = '  -1.92%  '

# Row 6
Set-TextValue $ws.Range("D6") '187.82'
$ws.Range("E6").Value = '  -3.65%  '

# Row 7
Set-TextValue $ws.Range("D7") '3.582.07'
$ws.Range("E7").Value = '  -1.28%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.622'
$ws.Range("E8").Value = '  -3.43%  '

# Row 9
Set-TextValue $ws.Range("D9") '1.00'
$ws.Range("E9").Value = '  +0.04%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.184'
$ws.Range("E10").Value = '  -0.49%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.654'
$ws.Range("E11").Value = '  -3.67%  '

# Row 12
Set-TextValue $ws.Range("D12") '55.40'
$ws.Range("E12").Value = '  -4.47%  '

# Row 13
Set-TextValue $ws.Range("D13") '0.0000306'
$ws.Range("E13").Value = '  -2.41%  '

# Row 14
Set-TextValue $ws.Range("D14") '9.58'
$ws.Range("E14").Value = '  -3.50%  '

# Row 15
Set-TextValue $ws.Range("D15") '4.162.45'
$ws.Range("E15").Value = '  -1.28%  '

# Row 16
Set-TextValue $ws.Range("D16") '19.71'
$ws.Range("E16").Value = '  -3.85%  '

# Row 17
Set-TextValue $ws.Range("D17") '3.576.53'
$ws.Range("E17").Value = '  -1.56%  '

# Row 18
Set-TextValue $ws.Range("D18") '70.160.58'
$ws.Range("E18").Value = '  -0.74%  '

# Row 19
Set-TextValue $ws.Range("D19") '12.61'
$ws.Range("E19").Value = '  -0.88%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.121'
$ws.Range("E20").Value = '  -1.08%  '

# Row 21
Set-TextValue $ws.Range("D21") '1.04'
$ws.Range("E21").Value = '  -2.82%  '

# Row 22
Set-TextValue $ws.Range("D22") '494.77'
$ws.Range("E22").Value = '  +1.31%  '

# Row 23
Set-TextValue $ws.Range("D23") '19.52'
$ws.Range("E23").Value = '  +0.99%  '

# Row 24
Set-TextValue $ws.Range("D24") '4.97'
$ws.Range("E24").Value = '  -5.44%  '

# Row 25
Set-TextValue $ws.Range("D25") '97.16'
$ws.Range("E25").Value = '  +6.38%  '

# Row 26
Set-TextValue $ws.Range("D26") '4.42'
$ws.Range("E26").Value = '  -1.79%  '

# Row 27
Set-TextValue $ws.Range("D27") '11.48'
$ws.Range("E27").Value = '  +0.52%  '

# Row 28
Set-TextValue $ws.Range("D28") '2.98'
$ws.Range("E28").Value = '  -6.03%  '

# Row 29
Set-TextValue $ws.Range("D29") '9.39'
$ws.Range("E29").Value = '  -1.81%  '

# Row 30
Set-TextValue $ws.Range("D30") '7.72'
$ws.Range("E30").Value = '  -2.87%  '

# Row 31
Set-TextValue $ws.Range("D31") '31.80'
$ws.Range("E31").Value = '  -3.09%  '

# Row 32
Set-TextValue $ws.Range("D32") '12.19'
$ws.Range("E32").Value = '  -0.87%  '

# Row 33
Set-TextValue $ws.Range("D33") '65.75'
$ws.Range("E33").Value = '  -0.85%  '

# Row 34
$ws.Range("E34").Value = '  -5.21%  '

# Row 35
Set-TextValue $ws.Range("D35") '573.30'
$ws.Range("E35").Value = '  -6.17%  '

# Row 36
Set-TextValue $ws.Range("D36") '3.21'
$ws.Range("E36").Value = '  +12.76%  '

# Row 37
Set-TextValue $ws.Range("D37") '39.16'
$ws.Range("E37").Value = '  -3.46%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.410'
$ws.Range("E38").Value = '  -0.20%  '

# Row 39
Set-TextValue $ws.Range("D39") '1.00'
$ws.Range("E39").Value = '  +0.21%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.0₃0799'
$ws.Range("E40").Value = '  -4.66%  '

# Row 41
$ws.Range("E41").Value = '  -2.23%  '

# Row 42
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D42") '3.19'
$ws.Range("E42").Value = '  +0.14%  '

# Row 43
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D43") '3.72'
$ws.Range("E43").Value = '  +11.48%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.135'
$ws.Range("E44").Value = '  -9.38%  '

# Row 45
Set-TextValue $ws.Range("D45") '3.09'
$ws.Range("E45").Value = '  -2.75%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.0457'
$ws.Range("E46").Value = '  -0.43%  '

# Row 47
Set-TextValue $ws.Range("D47") '3.203.19'
$ws.Range("E47").Value = '  -3.56%  '

# Row 48
Set-TextValue $ws.Range("D48") '9.55'
$ws.Range("E48").Value = '  -1.44%  '

# Row 49
$ws.Range("E49").Value = '  +31.08%  '

# Row 50
$ws.Range("E50").Value = '  -1.97%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.997'
$ws.Range("E51").Value = '  -0.23%  '
